$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Note: column D (Interruption) is set before columns B/C (Start/Stop Time)
# on each row so the dependent "Delta" formula in column E picks up the
# interruption minutes correctly when it recalculates.

# Row 58
$ws.Range("A58").Value = 41903
$ws.Range("D58").Value = 20
$ws.Range("B58").Value = 0.875
$ws.Range("C58").Value = 1.0138888888888888
$ws.Range("F58").Value = "Coding"

# Row 59
$ws.Range("A59").Value = 41904
$ws.Range("D59").Value = 10
$ws.Range("B59").Value = 0.53611111111111109
$ws.Range("C59").Value = 0.61319444444444449
$ws.Range("F59").Value = "Coding"

# Row 60
$ws.Range("A60").Value = 41904
$ws.Range("D60").Value = 15
$ws.Range("B60").Value = 0.98958333333333337
$ws.Range("C60").Value = 1.0770833333333334
$ws.Range("F60").Value = "Coding"

# Move selection to D61, matching the final state in the diff
$ws.Range("D61").Select()
